$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.038640994678112
$ws.Range("D2").Value = 1.046467931157318
$ws.Range("E2").Value = 1.037277884794242
$ws.Range("F2").Value = 1.05554291047988
$ws.Range("I2").Value = 1.044069507245682
$ws.Range("J2").Value = 1.043737331682886
$ws.Range("K2").Value = 1.04923328360805
$ws.Range("L2").Value = 1.040069226611013
$ws.Range("M2").Value = 1.058283080269713
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.039494528143958
$ws.Range("D3").Value = 1.047134865242459
$ws.Range("E3").Value = 1.038001261542631
$ws.Range("F3").Value = 1.056347585341504
$ws.Range("I3").Value = 1.044307145050531
$ws.Range("J3").Value = 1.044236222902161
$ws.Range("K3").Value = 1.049712242374992
$ws.Range("L3").Value = 1.040602606856677
$ws.Range("M3").Value = 1.058901239591474
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.040047410332944
$ws.Range("D4").Value = 1.047566897909925
$ws.Range("E4").Value = 1.038470216357446
$ws.Range("F4").Value = 1.0568690870278
$ws.Range("I4").Value = 1.044459910944364
$ws.Range("J4").Value = 1.044558964392123
$ws.Range("K4").Value = 1.050021958790645
$ws.Range("L4").Value = 1.040947950140412
$ws.Range("M4").Value = 1.059301396817266
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.04027998141337
$ws.Range("D5").Value = 1.047748638014806
$ws.Range("E5").Value = 1.038667574257091
$ws.Range("F5").Value = 1.057088521517427
$ws.Range("I5").Value = 1.044523893106046
$ws.Range("J5").Value = 1.044694625813416
$ws.Range("K5").Value = 1.050152113925172
$ws.Range("L5").Value = 1.041093181513329
$ws.Range("M5").Value = 1.059469661278389
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.04031903924992
$ws.Range("D6").Value = 1.047779159587002
$ws.Range("E6").Value = 1.038700723738194
$ws.Range("F6").Value = 1.057125376923871
$ws.Range("I6").Value = 1.044534621858596
$ws.Range("J6").Value = 1.044717402806023
$ws.Range("K6").Value = 1.050173964595776
$ws.Range("L6").Value = 1.041117569336995
$ws.Range("M6").Value = 1.059497915825776
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.040050517413311
$ws.Range("D7").Value = 1.047569325886047
$ws.Range("E7").Value = 1.038472852643214
$ws.Range("F7").Value = 1.056872018357656
$ws.Range("I7").Value = 1.044460766823637
$ws.Range("J7").Value = 1.044560777183318
$ws.Range("K7").Value = 1.050023698126634
$ws.Range("L7").Value = 1.040949890538076
$ws.Range("M7").Value = 1.059303645026461
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.038929327403362
$ws.Range("D8").Value = 1.046693223857837
$ws.Range("E8").Value = 1.037522169722768
$ws.Range("F8").Value = 1.055814682513331
$ws.Range("I8").Value = 1.044150024959647
$ws.Range("J8").Value = 1.043905948967054
$ws.Range("K8").Value = 1.049395190993966
$ws.Range("L8").Value = 1.040249440538164
$ws.Range("M8").Value = 1.058491954623183
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.036958228503736
$ws.Range("D9").Value = 1.045153186975931
$ws.Range("E9").Value = 1.035853776398905
$ws.Range("F9").Value = 1.053957910967526
$ws.Range("I9").Value = 1.043594822964747
$ws.Range("J9").Value = 1.042751543562724
$ws.Range("K9").Value = 1.048286192167221
$ws.Range("L9").Value = 1.039016831250581
$ws.Range("M9").Value = 1.057062999444251
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.035647337271497
$ws.Range("D10").Value = 1.0441291327737
$ws.Range("E10").Value = 1.03474621353605
$ws.Range("F10").Value = 1.052724472659743
$ws.Range("I10").Value = 1.043219603827966
$ws.Range("J10").Value = 1.041981669605892
$ws.Range("K10").Value = 1.047545936146493
$ws.Range("L10").Value = 1.038196297730836
$ws.Range("M10").Value = 1.056111361924702
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.035080479736297
$ws.Range("D11").Value = 1.043686353520094
$ws.Range("E11").Value = 1.034267763660728
$ws.Range("F11").Value = 1.052191450069494
$ws.Range("I11").Value = 1.043055934115908
$ws.Range("J11").Value = 1.041648257454574
$ws.Range("K11").Value = 1.04722519443128
$ws.Range("L11").Value = 1.037841300805957
$ws.Range("M11").Value = 1.055699549115976
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.034870040429723
$ws.Range("D12").Value = 1.043521983972395
$ws.Range("E12").Value = 1.034090217987622
$ws.Range("F12").Value = 1.051993623589091
$ws.Range("I12").Value = 1.042994960819238
$ws.Range("J12").Value = 1.041524406850413
$ws.Range("K12").Value = 1.047106027042032
$ws.Range("L12").Value = 1.037709485516043
$ws.Range("M12").Value = 1.055546623027903
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.034915175063435
$ws.Range("D13").Value = 1.043557237317762
$ws.Range("E13").Value = 1.034128294328538
$ws.Range("F13").Value = 1.052036050685858
$ws.Range("I13").Value = 1.043008047888995
$ws.Range("J13").Value = 1.041550973498482
$ws.Range("K13").Value = 1.047131590171572
$ws.Range("L13").Value = 1.037737758222551
$ws.Range("M13").Value = 1.055579424378842
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.035063082365266
$ws.Range("D14").Value = 1.043672764664746
$ws.Range("E14").Value = 1.034253084162497
$ws.Range("F14").Value = 1.05217509435817
$ws.Range("I14").Value = 1.043050897692634
$ws.Range("J14").Value = 1.041638020049944
$ws.Range("K14").Value = 1.047215344614683
$ws.Range("L14").Value = 1.03783040395899
$ws.Range("M14").Value = 1.055686907382352
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.035154228462925
$ws.Range("D15").Value = 1.04374395796593
$ws.Range("E15").Value = 1.034329994132415
$ws.Range("F15").Value = 1.052260785257347
$ws.Range("I15").Value = 1.043077275165538
$ws.Range("J15").Value = 1.04169165148184
$ws.Range("K15").Value = 1.047266944603561
$ws.Range("L15").Value = 1.037887492241121
$ws.Range("M15").Value = 1.055753136493127
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.035684974034667
$ws.Range("D16").Value = 1.044158532273405
$ws.Range("E16").Value = 1.034777990687135
$ws.Range("F16").Value = 1.05275987020338
$ws.Range("I16").Value = 1.043230440890043
$ws.Range("J16").Value = 1.042003796092852
$ws.Range("K16").Value = 1.047567218489994
$ws.Range("L16").Value = 1.038219864135357
$ws.Range("M16").Value = 1.05613869804028
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.036018103180594
$ws.Range("D17").Value = 1.044418757184731
$ws.Range("E17").Value = 1.03505931148521
$ws.Range("F17").Value = 1.053073219325814
$ws.Range("I17").Value = 1.043326197709186
$ws.Range("J17").Value = 1.042199583231984
$ws.Range("K17").Value = 1.0477555181457
$ws.Range("L17").Value = 1.038428433401537
$ws.Range("M17").Value = 1.056380619338975
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03621248567502
$ws.Range("D18").Value = 1.044570603891076
$ws.Range("E18").Value = 1.035223510261544
$ws.Range("F18").Value = 1.053256093130611
$ws.Range("I18").Value = 1.043381935460405
$ws.Range("J18").Value = 1.042313777489123
$ws.Range("K18").Value = 1.047865330187814
$ws.Range("L18").Value = 1.038550117131249
$ws.Range("M18").Value = 1.056521752197331
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.036278777568811
$ws.Range("D19").Value = 1.04462239014202
$ws.Range("E19").Value = 1.03527951625328
$ws.Range("F19").Value = 1.053318465709677
$ws.Range("I19").Value = 1.043400920963964
$ws.Range("J19").Value = 1.04235271388129
$ws.Range("K19").Value = 1.047902769830743
$ws.Range("L19").Value = 1.038591612982489
$ws.Range("M19").Value = 1.056569878921084
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.0359823539207
$ws.Range("D20").Value = 1.044390831097134
$ws.Range("E20").Value = 1.035029117126621
$ws.Range("F20").Value = 1.053039589307914
$ws.Range("I20").Value = 1.043315935857998
$ws.Range("J20").Value = 1.042178577639949
$ws.Range("K20").Value = 1.047735317429927
$ws.Range("L20").Value = 1.038406052903733
$ws.Range("M20").Value = 1.056354660952591
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.035019524113781
$ws.Range("D21").Value = 1.04363874201038
$ws.Range("E21").Value = 1.034216331906437
$ws.Range("F21").Value = 1.052134144979046
$ws.Range("I21").Value = 1.043038284430387
$ws.Range("J21").Value = 1.041612387178934
$ws.Range("K21").Value = 1.047190681826268
$ws.Range("L21").Value = 1.037803120807532
$ws.Range("M21").Value = 1.055655255212933
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.034414831479086
$ws.Range("D22").Value = 1.04316644383893
$ws.Range("E22").Value = 1.033706297117698
$ws.Range("F22").Value = 1.051565793521597
$ws.Range("I22").Value = 1.042862678268857
$ws.Range("J22").Value = 1.041256363868723
$ws.Range("K22").Value = 1.046848076994051
$ws.Range("L22").Value = 1.037424302456547
$ws.Range("M22").Value = 1.055215740530317
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.0347353258142
$ws.Range("D23").Value = 1.043416763403466
$ws.Range("E23").Value = 1.03397658123651
$ws.Range("F23").Value = 1.051866997887931
$ws.Range("I23").Value = 1.042955868318053
$ws.Range("J23").Value = 1.041445101633601
$ws.Range("K23").Value = 1.047029714051478
$ws.Range("L23").Value = 1.037625095275227
$ws.Range("M23").Value = 1.055448713298895
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.035998507244282
$ws.Range("D24").Value = 1.044403449502964
$ws.Range("E24").Value = 1.035042760320654
$ws.Range("F24").Value = 1.053054784948757
$ws.Range("I24").Value = 1.043320573104055
$ws.Range("J24").Value = 1.04218806917863
$ws.Range("K24").Value = 1.047744445325692
$ws.Range("L24").Value = 1.0384161655985
$ws.Range("M24").Value = 1.056366390355275
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.037467252191233
$ws.Range("D25").Value = 1.045550866266181
$ws.Range("E25").Value = 1.036284275545517
$ws.Range("F25").Value = 1.054437162046645
$ws.Range("I25").Value = 1.043739255752403
$ws.Range("J25").Value = 1.043050038212056
$ws.Range("K25").Value = 1.04857306339714
$ws.Range("L25").Value = 1.0393352836661
$ws.Range("M25").Value = 1.057432249334854

Write-Output "Updated 240 cells"
